$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "Mehr darauf achten Zeit nicht zu unterschätzen!"
$ws.Range("A8").Value = "Effort von allen 3 Mitgliedern beachten!"
$ws.Range("A9").Value = "Bessere Aufteilung der Aufgaben"
$ws.Range("A6").Value = "Notes:"

$ws.Range("A10").Select()
